$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (existing rows 13-15 shift down to 14-16),
# inheriting the formatting (style) of the row that was at 13.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the "Possible_Problem" answer for the
# "electrical accessories" question (A15/A16 after the shift), reusing
# the same Node/Relationship/Node2 content already used in row 7.
$ws.Range("A13").Value = $ws.Range("A14").Value()
$ws.Range("B13").Value = $ws.Range("B7").Value()
$ws.Range("C13").Value = $ws.Range("C7").Value()

# Match the row height used for the equivalent "Possible_Problem" row (row 7).
$ws.Rows.Item(13).RowHeight = 345.6

# Update the view so the newly inserted row's C cell is selected/scrolled into view.
$ws.Range("C13").Select() | Out-Null
